$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLoginTest")
$ws.Range("A4").Value = "Kim"
$ws.Range("B4").Value = "kim124"
$ws.Range("C4").Value = "Invalid credentials"
